$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: code 101, ФИЗИКА, ТАРБИЯИ ҶИСМОНӢ, 8:00-8:50, БОЙМАТОВ Д., Пониделник
$ws.Range("A2").Value = 101
$ws.Range("B2").Value = "ФИЗИКА "
$ws.Range("C2").Value = "ТАРБИЯИ ҶИСМОНӢ"
$ws.Range("D2").Value = "8:00-8:50"
$ws.Range("E2").Value = "БОЙМАТОВ Д."
$ws.Range("F2").Value = "Пониделник"

# Row 3: code 101, ФИЗИКА, ТАРБИЯИ ҶИСМОНӢ, 8:00-8:50, БОЙМАТОВ Д., Вторник
$ws.Range("A3").Value = 101
$ws.Range("B3").Value = "ФИЗИКА "
$ws.Range("C3").Value = "ТАРБИЯИ ҶИСМОНӢ"
$ws.Range("D3").Value = "8:00-8:50"
$ws.Range("E3").Value = "БОЙМАТОВ Д."
$ws.Range("F3").Value = "Вторник"

# Row 4: code 102, ФИЗИКА- МАТЕМАТИКА, ТАРБИЯИ ҶИСМОНӢ, 8:00-8:50, БОЙМАТОВ Д., Пониделник
$ws.Range("A4").Value = 102
$ws.Range("B4").Value = "ФИЗИКА- МАТЕМАТИКА"
$ws.Range("C4").Value = "ТАРБИЯИ ҶИСМОНӢ"
$ws.Range("D4").Value = "8:00-8:50"
$ws.Range("E4").Value = "БОЙМАТОВ Д."
$ws.Range("F4").Value = "Пониделник"

# Row 5: code 102, ФИЗИКА- МАТЕМАТИКА, ТАРБИЯИ ҶИСМОНӢ, 9:00-10:50, БОЙМАТОВ Д., Пониделник
$ws.Range("A5").Value = 102
$ws.Range("B5").Value = "ФИЗИКА- МАТЕМАТИКА"
$ws.Range("C5").Value = "ТАРБИЯИ ҶИСМОНӢ"
$ws.Range("D5").Value = "9:00-10:50"
$ws.Range("E5").Value = "БОЙМАТОВ Д."
$ws.Range("F5").Value = "Пониделник"

# Row 6: code 102, ФИЗИКА- МАТЕМАТИКА, ТАРБИЯИ ҶИСМОНӢ, 10:00-11:50, БОЙМАТОВ Д., Пониделник
$ws.Range("A6").Value = 102
$ws.Range("B6").Value = "ФИЗИКА- МАТЕМАТИКА"
$ws.Range("C6").Value = "ТАРБИЯИ ҶИСМОНӢ"
$ws.Range("D6").Value = "10:00-11:50"
$ws.Range("E6").Value = "БОЙМАТОВ Д."
$ws.Range("F6").Value = "Пониделник"

# Row 7: code 102, ФИЗИКА- МАТЕМАТИКА, ТАРБИЯИ ҶИСМОНӢ, 11:00-11:50, БОЙМАТОВ Д., Пониделник
$ws.Range("A7").Value = 102
$ws.Range("B7").Value = "ФИЗИКА- МАТЕМАТИКА"
$ws.Range("C7").Value = "ТАРБИЯИ ҶИСМОНӢ"
$ws.Range("D7").Value = "11:00-11:50"
$ws.Range("E7").Value = "БОЙМАТОВ Д."
$ws.Range("F7").Value = "Пониделник"

# Row 8: code 103, ФИЗИКА- Информатика,  ҶИСМОНӢ, 8:00-11:50, БОЙМАТОВ Д., Пониделник
$ws.Range("A8").Value = 103
$ws.Range("B8").Value = "ФИЗИКА- Информатика"
$ws.Range("C8").Value = " ҶИСМОНӢ"
$ws.Range("D8").Value = "8:00-11:50"
$ws.Range("E8").Value = "БОЙМАТОВ Д."
$ws.Range("F8").Value = "Пониделник"

# Row 9: code 201, Физика, ТАРБИЯИ , 9:00-11:50, БОЙМАТОВ Д., Пониделник
$ws.Range("A9").Value = 201
$ws.Range("B9").Value = "Физика"
$ws.Range("C9").Value = "ТАРБИЯИ "
$ws.Range("D9").Value = "9:00-11:50"
$ws.Range("E9").Value = "БОЙМАТОВ Д."
$ws.Range("F9").Value = "Пониделник"

# Column widths (COM ColumnWidth snaps to 1/6-char + 5/6 padding internally,
# so we pre-compensate to land as close as possible to the target stored width)
$ws.Columns.Item(2).ColumnWidth = 27.166666666666668
$ws.Columns.Item(3).ColumnWidth = 24.498697916666668
$ws.Columns.Item(5).ColumnWidth = 12.276041666666666
$ws.Columns.Item(6).ColumnWidth = 12.276041666666666

# View / zoom / selection
$excel.ActiveWindow.Zoom = 120
$ws.Range("C9").Select() | Out-Null
